$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update variable-dictionary cell values (new clinical-trial "other" device columns) ---
# (written in this order so new shared-string entries land in the same order
# the target workbook has them in)

# Row 11 ("other" section): device_dt_tm -> TIMESTAMP, width 10 -> 14
$ws.Range("A11").Value = "TIMESTAMP"
$ws.Range("D11").Value = 14

# Row 12 ("other" section): value -> GLUCOSE, width 10 -> 14, and it becomes a
# header-styled row like rows 7-9 (copy that style before overwriting the value)
$ws.Range("A7").Copy($ws.Range("A12"))
$ws.Range("A12").Value = "GLUCOSE"
$ws.Range("D12").Value = 14

# Row 10 ("other" section): pt_id -> ID_VISIT_DEVICEID, width 10 -> 14
$ws.Range("A10").Value = "ID_VISIT_DEVICEID"
$ws.Range("D10").Value = 14

# Row 2: libre id -> deviceid
$ws.Range("B2").Value = "deviceid"

# Row 9: dexcomg6 device -> deviceid
$ws.Range("B9").Value = "deviceid"

# Row 1 header: sensor_lifetime -> expecteddaysofwear
$ws.Range("D1").Value = "expecteddaysofwear"

# --- View state: zoom to 153% and move the selection to F9 ---
$excel.ActiveWindow.Zoom = 153
$ws.Range("F9").Select()
